# Re-point the table on slide 6 at the new (built-in) table style.
# (was {21F73A81-A468-4E02-8142-FFCCF8D3A25F} -> now {FFD72DC7-0C5E-42CF-B94A-D327153BB96E})
$p = $ppt.ActivePresentation
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{FFD72DC7-0C5E-42CF-B94A-D327153BB96E}")

# Swap the deck's theme palette from "Integral" to the default "Office Theme"
# colours (the other half of the edit - the design/theme was changed).
# ThemeColorScheme indices follow the standard clrScheme order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$cs = $p.Slides.Item(1).ThemeColorScheme
$cs.Item(1).RGB  = 0         # dk1      000000
$cs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      44546A
$cs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  FFC000
$cs.Item(9).RGB  = 12874308  # accent5  4472C4
$cs.Item(10).RGB = 4697456   # accent6  70AD47
$cs.Item(11).RGB = 12673797  # hlink    0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72
